{"js": "const replacements = [\n    [\"906\u00d76=\", \"236\u00d79=\"],\n    [\"984\u00d76=\", \"781\u00d76=\"],\n    [\"381\u00d73=\", \"470\u00d79=\"],\n    [\"337\u00d77=\", \"829\u00d73=\"],\n    [\"146\u00d75=\", \"925\u00d78=\"],\n    [\"682\u00d73=\", \"540\u00d78=\"],\n    [\"229\u00d73=\", \"824\u00d76=\"],\n    [\"342\u00d76=\", \"897\u00d79=\"],\n    [\"200\u00d78=\", \"106\u00d75=\"],\n    [\"884\u00d79=\", \"215\u00d76=\"],\n    [\"860\u00d77=\", \"253\u00d74=\"],\n    [\"272\u00d79=\", \"515\u00d79=\"],\n    [\"984\u00d73=\", \"290\u00d78=\"],\n    [\"978\u00d74=\", \"612\u00d76=\"],\n    [\"977\u00d75=\", \"870\u00d77=\"],\n    [\"832\u00d75=\", \"597\u00d75=\"],\n    [\"517\u00d75=\", \"450\u00d78=\"],\n    [\"936\u00d75=\", \"651\u00d74=\"],\n    [\"489\u00d77=\", \"611\u00d77=\"],\n    [\"917\u00d72=\", \"131\u00d74=\"],\n    [\"265\u00d74=\", \"974\u00d77=\"],\n    [\"470\u00d75=\", \"486\u00d72=\"],\n    [\"386\u00d75=\", \"884\u00d78=\"],\n    [\"270\u00d77=\", \"632\u00d76=\"],\n    [\"169\u00d72=\", \"841\u00d74=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n    const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n    results.load(\"text\");\n    await context.sync();\n\n    for (let i = 0; i < results.items.length; i++) {\n        results.items[i].insertText(newText, \"Replace\");\n    }\n    await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n    @{ Old = \"906\u00d76=\"; New = \"236\u00d79=\" }\n    @{ Old = \"984\u00d76=\"; New = \"781\u00d76=\" }\n    @{ Old = \"381\u00d73=\"; New = \"470\u00d79=\" }\n    @{ Old = \"337\u00d77=\"; New = \"829\u00d73=\" }\n    @{ Old = \"146\u00d75=\"; New = \"925\u00d78=\" }\n    @{ Old = \"682\u00d73=\"; New = \"540\u00d78=\" }\n    @{ Old = \"229\u00d73=\"; New = \"824\u00d76=\" }\n    @{ Old = \"342\u00d76=\"; New = \"897\u00d79=\" }\n    @{ Old = \"200\u00d78=\"; New = \"106\u00d75=\" }\n    @{ Old = \"884\u00d79=\"; New = \"215\u00d76=\" }\n    @{ Old = \"860\u00d77=\"; New = \"253\u00d74=\" }\n    @{ Old = \"272\u00d79=\"; New = \"515\u00d79=\" }\n    @{ Old = \"984\u00d73=\"; New = \"290\u00d78=\" }\n    @{ Old = \"978\u00d74=\"; New = \"612\u00d76=\" }\n    @{ Old = \"977\u00d75=\"; New = \"870\u00d77=\" }\n    @{ Old = \"832\u00d75=\"; New = \"597\u00d75=\" }\n    @{ Old = \"517\u00d75=\"; New = \"450\u00d78=\" }\n    @{ Old = \"936\u00d75=\"; New = \"651\u00d74=\" }\n    @{ Old = \"489\u00d77=\"; New = \"611\u00d77=\" }\n    @{ Old = \"917\u00d72=\"; New = \"131\u00d74=\" }\n    @{ Old = \"265\u00d74=\"; New = \"974\u00d77=\" }\n    @{ Old = \"470\u00d75=\"; New = \"486\u00d72=\" }\n    @{ Old = \"386\u00d75=\"; New = \"884\u00d78=\" }\n    @{ Old = \"270\u00d77=\"; New = \"632\u00d76=\" }\n    @{ Old = \"169\u00d72=\"; New = \"841\u00d74=\" }\n)\n\nforeach ($pair in $pairs) {\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Text = $pair.Old\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $pair.New\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($find.Text, $find.MatchCase, $find.MatchWholeWord, $find.MatchWildcards, $false, $false, $find.Forward, $find.Wrap, $false, $find.Replacement.Text, 2) | Out-Null\n}\n"}
